$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 10.41083696637312
$ws.Range("E2").Value = 10000
$ws.Range("F2").Value = 0.0000000009025885347604574
$ws.Range("D3").Value = 10.41064597649401
$ws.Range("E3").Value = 6652
$ws.Range("F3").Value = 0.00000000007061195356961412
$ws.Range("D4").Value = 10.40660457794342
$ws.Range("E4").Value = 182
$ws.Range("F4").Value = 0.00000000006732421208540055
$ws.Range("D5").Value = 10.40660333716813
$ws.Range("E5").Value = 176
$ws.Range("F5").Value = 0.00000000006938399862840957
$ws.Range("D8").Value = 10.40660457794342
$ws.Range("E8").Value = 182
$ws.Range("F8").Value = 0.00000000006782656788632418
$ws.Range("D9").Value = 10.40660291773662
$ws.Range("E9").Value = 174
$ws.Range("F9").Value = 0.00000000007387721035772134
$ws.Range("D10").Value = 199.9999999956344
$ws.Range("E10").Value = 18
$ws.Range("F10").Value = 0.00000000006548376063762736
$ws.Range("D11").Value = 199.9999284744817
$ws.Range("E11").Value = 11
$ws.Range("F11").Value = 0.00000000006650680006685051
$ws.Range("D12").Value = 199.9999999999829
$ws.Range("E12").Value = 22
$ws.Range("F12").Value = 0.00000000003517286017872089
$ws.Range("D13").Value = 199.9999999999829
$ws.Range("E13").Value = 22
$ws.Range("F13").Value = 0.00000000003283560090387037
$ws.Range("D16").Value = 199.9999999999829
$ws.Range("E16").Value = 22
$ws.Range("F16").Value = 0.00000000003502933054610741
$ws.Range("D17").Value = 199.9999999999829
$ws.Range("E17").Value = 22
$ws.Range("F17").Value = 0.00000000003550752580976184
$ws.Range("D18").Value = 6197287.051397814
$ws.Range("E18").Value = 119
$ws.Range("F18").Value = 0.00000000009691242166203807
$ws.Range("D19").Value = 6197285.757615398
$ws.Range("E19").Value = 75
$ws.Range("F19").Value = 0.0000000000954792100668907
$ws.Range("D20").Value = 6197287.048449507
$ws.Range("E20").Value = 115
$ws.Range("F20").Value = 0.00000000009402571091828995
$ws.Range("D21").Value = 6197287.047441144
$ws.Range("E21").Value = 114
$ws.Range("F21").Value = 0.00000000006832363075419507
$ws.Range("D24").Value = 6197287.048449507
$ws.Range("E24").Value = 115
$ws.Range("F24").Value = 0.00000000009378466327958547
$ws.Range("D25").Value = 6197286.989795357
$ws.Range("E25").Value = 98
$ws.Range("F25").Value = 0.00000000008360761474517046
$ws.Range("D26").Value = 30148.79411856191
$ws.Range("E26").Value = 1231
$ws.Range("F26").Value = 0.00000000009906775418343501
$ws.Range("D27").Value = 30148.6402523005
$ws.Range("E27").Value = 595
$ws.Range("F27").Value = 0.00000000009654664719054071
$ws.Range("D28").Value = 30148.78281109525
$ws.Range("E28").Value = 859
$ws.Range("F28").Value = 0.00000000009609852570279945
$ws.Range("D29").Value = 30135.17405261119
$ws.Range("E29").Value = 128
$ws.Range("F29").Value = 0.00000000009695341624095958
$ws.Range("D32").Value = 30148.75084112083
$ws.Range("E32").Value = 724
$ws.Range("F32").Value = 0.00000000006670350626795326
$ws.Range("D33").Value = 30146.64990782431
$ws.Range("E33").Value = 325
$ws.Range("F33").Value = 0.00000000008378623029590728
$ws.Range("D34").Value = 5.502378378978929
$ws.Range("E34").Value = 19
$ws.Range("F34").Value = 0.00000000004614223178292824
$ws.Range("D35").Value = 5.50237842977301
$ws.Range("E35").Value = 14
$ws.Range("F35").Value = 0.00000000006065363656690676
$ws.Range("D36").Value = 5.502378378877919
$ws.Range("E36").Value = 22
$ws.Range("F36").Value = 0.00000000008892165939451615
$ws.Range("D37").Value = 5.502378378877919
$ws.Range("E37").Value = 22
$ws.Range("F37").Value = 0.00000000008274454662063049
$ws.Range("D40").Value = 5.502378378877919
$ws.Range("E40").Value = 22
$ws.Range("F40").Value = 0.00000000008891972238831549
$ws.Range("D41").Value = 5.502378378877919
$ws.Range("E41").Value = 22
$ws.Range("F41").Value = 0.00000000008328448708518714
$ws.Range("D42").Value = 0.9788022306513534
$ws.Range("E42").Value = 4875
$ws.Range("F42").Value = 0.00000000009991530743457809
$ws.Range("D43").Value = 0.9787919201260611
$ws.Range("E43").Value = 2234
$ws.Range("F43").Value = 0.00000000009630889006611835
$ws.Range("D44").Value = 0.9787705852933443
$ws.Range("E44").Value = 1675
$ws.Range("F44").Value = 0.00000000007899018222984142
$ws.Range("D45").Value = 0.9785554077148678
$ws.Range("E45").Value = 545
$ws.Range("F45").Value = 0.00000000006763274964508732
$ws.Range("D48").Value = 0.9787950509179197
$ws.Range("E48").Value = 2412
$ws.Range("F48").Value = 0.0000000000953195032114241
$ws.Range("D49").Value = 0.9784064008082318
$ws.Range("E49").Value = 205
$ws.Range("F49").Value = 0.00000000008052059947469192
$ws.Range("D50").Value = 3015179080.010438
$ws.Range("E50").Value = 520
$ws.Range("F50").Value = 0.0000000000995572108325373
$ws.Range("D51").Value = 3015164921.695023
$ws.Range("E51").Value = 277
$ws.Range("F51").Value = 0.00000000009801165388238914
$ws.Range("D52").Value = 3015178990.986972
$ws.Range("E52").Value = 443
$ws.Range("F52").Value = 0.00000000009596636076646659
$ws.Range("D53").Value = 3015178621.431344
$ws.Range("E53").Value = 391
$ws.Range("F53").Value = 0.000000000096114871223921
$ws.Range("D56").Value = 3015178893.111733
$ws.Range("E56").Value = 420
$ws.Range("F56").Value = 0.0000000000865302926825525
$ws.Range("D57").Value = 3006240323.348099
$ws.Range("E57").Value = 54
$ws.Range("F57").Value = 0.00000000001143794261840101
$ws.Range("D58").Value = 5.502378378573628
$ws.Range("E58").Value = 19
$ws.Range("F58").Value = 0.00000000004614239320351041
$ws.Range("D59").Value = 5.502378429367709
$ws.Range("E59").Value = 14
$ws.Range("F59").Value = 0.00000000006063232950316975
$ws.Range("D60").Value = 5.502378378472619
$ws.Range("E60").Value = 22
$ws.Range("F60").Value = 0.0000000000888863090379013
$ws.Range("D61").Value = 5.502378378472619
$ws.Range("E61").Value = 22
$ws.Range("F61").Value = 0.00000000008274131828305741
$ws.Range("D64").Value = 5.502378378472619
$ws.Range("E64").Value = 22
$ws.Range("F64").Value = 0.00000000008895684834704741
$ws.Range("D65").Value = 5.502378378472619
$ws.Range("E65").Value = 22
$ws.Range("F65").Value = 0.00000000008249531848194686
$ws.Range("D66").Value = 66.4965639778705
$ws.Range("E66").Value = 271
$ws.Range("F66").Value = 0.00000000009581047272545482
$ws.Range("D67").Value = 66.49648730538858
$ws.Range("E67").Value = 170
$ws.Range("F67").Value = 0.00000000009989956374508951
$ws.Range("D68").Value = 66.49656370279384
$ws.Range("E68").Value = 250
$ws.Range("F68").Value = 0.00000000009681297783983837
$ws.Range("D69").Value = 66.49656326260559
$ws.Range("E69").Value = 238
$ws.Range("F69").Value = 0.00000000009477078382019032
$ws.Range("D72").Value = 66.49656367793219
$ws.Range("E72").Value = 249
$ws.Range("F72").Value = 0.00000000009422369045765888
$ws.Range("D73").Value = 66.49655852666243
$ws.Range("E73").Value = 209
$ws.Range("F73").Value = 0.00000000008822362850942354
$ws.Range("D74").Value = 11.86733830706314
$ws.Range("E74").Value = 870
$ws.Range("F74").Value = 0.00000000009887665439129894
$ws.Range("D75").Value = 11.86722648600631
$ws.Range("E75").Value = 329
$ws.Range("F75").Value = 0.00000000009935833806553928
$ws.Range("D76").Value = 11.86733623746636
$ws.Range("E76").Value = 628
$ws.Range("F76").Value = 0.00000000007899347461762838
$ws.Range("D77").Value = 11.86732786565626
$ws.Range("E77").Value = 508
$ws.Range("F77").Value = 0.00000000003536210088174771
$ws.Range("D80").Value = 11.86733586681646
$ws.Range("E80").Value = 616
$ws.Range("F80").Value = 0.00000000007083267537190448
$ws.Range("D81").Value = 11.86710093323266
$ws.Range("E81").Value = 272
$ws.Range("F81").Value = 0.00000000007929097315067625
$ws.Range("D82").Value = 12.60584912893822
$ws.Range("E82").Value = 2832
$ws.Range("F82").Value = 0.00000000004536838100943374
$ws.Range("D83").Value = 12.60584912953977
$ws.Range("E83").Value = 2833
$ws.Range("F83").Value = 0.00000000004519294147690548
$ws.Range("D84").Value = 12.60614935431771
$ws.Range("E84").Value = 1957
$ws.Range("F84").Value = 0.000000000001919079579004
$ws.Range("D85").Value = 12.60637048546948
$ws.Range("E85").Value = 1893
$ws.Range("F85").Value = 0.00000000004204879152795875
$ws.Range("D88").Value = 12.60932751697525
$ws.Range("E88").Value = 1304
$ws.Range("F88").Value = 0.00000000008831218358091898
$ws.Range("D89").Value = 12.60185161443958
$ws.Range("E89").Value = 719
$ws.Range("F89").Value = 0.00000000002779449023045671
$ws.Range("D90").Value = 4359568.100322182
$ws.Range("E90").Value = 92
$ws.Range("F90").Value = 0.0000000000918432711756057
$ws.Range("D91").Value = 4359559.206934161
$ws.Range("E91").Value = 56
$ws.Range("F91").Value = 0.00000000009542921829765347
$ws.Range("D92").Value = 4359568.100636504
$ws.Range("E92").Value = 93
$ws.Range("F92").Value = 0.00000000009520576418355374
$ws.Range("D93").Value = 4359568.100636504
$ws.Range("E93").Value = 93
$ws.Range("F93").Value = 0.00000000009028101503111049
$ws.Range("D96").Value = 4359568.100636504
$ws.Range("E96").Value = 93
$ws.Range("F96").Value = 0.0000000000973292191167486
$ws.Range("D97").Value = 4359568.096880184
$ws.Range("E97").Value = 87
$ws.Range("F97").Value = 0.00000000008161009838585245
